# Daily attendance processing - 2026-01-24 11:33:47
# Normalises the "Recorded By" (column G) value ordering for rows where the
# "System" (or "system") entry was listed first, moving it after the other
# recorder(s) instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Exact old-value -> new-value replacements observed for the "Recorded By" column.
$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
